$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data corrections
# ---------------------------------------------------------------------------

# "User Stories" sheet (Sheet2 codeName) - self-assessment score swap between
# two students (US corrections)
$wsUserStories = $wb.Worksheets.Item("User Stories")
$wsUserStories.Range("C6").Value = 5
$wsUserStories.Range("C10").Value = 4

# "Project Development" sheet (Sheet4 codeName) - small correction, raising
# two evaluation scores from 4 to 5
$wsProjectDev = $wb.Worksheets.Item("Project Development")
$wsProjectDev.Range("C4").Value = 5
$wsProjectDev.Range("E4").Value = 5
$wsProjectDev.Range("C5").Value = 5
$wsProjectDev.Range("E5").Value = 5

# "Project Management" sheet (Sheet5 codeName) - small correction in US04,
# raising evaluation scores
$wsProjectMgmt = $wb.Worksheets.Item("Project Management")
$wsProjectMgmt.Range("E7").Value = 5
$wsProjectMgmt.Range("F7").Value = 4
$wsProjectMgmt.Range("C8").Value = 5
$wsProjectMgmt.Range("D8").Value = 5
$wsProjectMgmt.Range("E8").Value = 5
$wsProjectMgmt.Range("F8").Value = 5

# ---------------------------------------------------------------------------
# View / selection state (mirrors the author re-navigating the workbook
# before saving - continuing work by copying sprint B into the next sprint)
# ---------------------------------------------------------------------------

$wsGroupSelf = $wb.Worksheets.Item("Group and Self Assessment")
$wsGroupSelf.Activate()
$wsGroupSelf.Range("V10").Select()

$wsUserStories.Activate()
$excel.ActiveWindow.Zoom = 72
$wsUserStories.Range("C11").Select()

$wsProjectDev.Activate()
$wsProjectDev.Range("C5").Select()

# Leave "Project Management" as the active/selected sheet and cell, matching
# the saved workbook state (activeTab + tabSelected).
$wsProjectMgmt.Activate()
$wsProjectMgmt.Range("E5").Select()
